$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ORD-2413-8980-9419"
$ws.Range("H2").Value = "ORD-2413-8980-9419"
$ws.Range("F2").Value = "SUB-0356-5642-8669"

$ws.Range("F2").Select()
